# ---------------------------------------------------------------------------
# Applies the "edit window" fix-up described by the commit:
#   "vylepsovani okna editu - jeste chyby pri ukladani zmen"
# (improving the edit window - still bugs when saving changes)
#
# Net effect on the workbook:
#   * Sheet "ip_address_list"   : rows 4-8 get new content, rows 9-12 removed
#   * Sheet "ip_adress_fav_list": rows 1-3 get new content, row 4 removed
#   * Sheet "Settings"          : B1 4->0, B3 1->0
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1) ip_address_list
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ip_address_list")

# Row 4: 511_Teleflex -> 514_Teleflex
$ws1.Cells.Item(4,1).Value = "514_Teleflex"
$ws1.Cells.Item(4,2).Value = "192.168.14.240"
$ws1.Cells.Item(4,3).Value = "255.255.255.0"
$ws1.Cells.Item(4,4).Value = "PC:192.168.14.240"
$ws1.Cells.Item(4,5).Value = $true

# Row 5: 514_Teleflex -> 518_Valeo (notes cleared)
$ws1.Cells.Item(5,1).Value = "518_Valeo"
$ws1.Cells.Item(5,2).Value = "192.168.208.242"
$ws1.Cells.Item(5,3).Value = "255.255.255.0"
$ws1.Cells.Item(5,4).ClearContents()
$ws1.Cells.Item(5,5).Value = 0

# Row 6: 518_Valeo -> 527_Teijin
$ws1.Cells.Item(6,1).Value = "527_Teijin"
$ws1.Cells.Item(6,2).Value = "10.101.28.176"
$ws1.Cells.Item(6,3).Value = "255.255.255.0"
$ws1.Cells.Item(6,4).Value = "XG-X2900:`t`t10.101.28.175"
$ws1.Cells.Item(6,5).Value = $true

# Row 7: 529_Witte -> 515_ZF Stara Boleslav
$ws1.Cells.Item(7,1).Value = "515_ZF Stara Boleslav"
$ws1.Cells.Item(7,2).Value = "10.9.250.240"
$ws1.Cells.Item(7,3).Value = "255.255.255.0"
$ws1.Cells.Item(7,4).Value = "NAS - 10.9.250.100`nUser:spravce Pass:Jhv*2708 `nUser:jhvadmin Pass:jhvadm1n >>> na portu 8080. `n123TPV456"
$ws1.Cells.Item(7,5).Value = 0

# Row 8: Domaci Wifi -> 503_Witte
$ws1.Cells.Item(8,1).Value = "503_Witte"
$ws1.Cells.Item(8,2).Value = "192.168.0.240"
$ws1.Cells.Item(8,3).Value = "255.255.255.0"
$ws1.Cells.Item(8,4).Value = "PC:`t10.96.205.175`nNAS:`t10.96.205.166`nFH:`t10.96.205.154`n`t10.96.205.267`n-----------------------------------------`nuser:JHV_Vision, omron `nPass:*Jhv2708`n---------------------------------------`nFortiClient Austin: `nPass:`n1Pm#J@PFIkzM&Q@i `nUVt1@Ex2p78kxp30atD7we@!qGK"
$ws1.Cells.Item(8,5).Value = 1

# Rows 9-12 no longer exist (old favourites / extra entries folded away)
$ws1.Range("A9:E12").EntireRow.Delete()

# ----------------------------------------------------------------------
# 2) ip_adress_fav_list  (favourites mirror the flagged rows above)
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ip_adress_fav_list")

# Row 1: 518_Valeo II -> 514_Teleflex
$ws2.Cells.Item(1,1).Value = "514_Teleflex"
$ws2.Cells.Item(1,2).Value = "192.168.14.240"
$ws2.Cells.Item(1,3).Value = "255.255.255.0"
$ws2.Cells.Item(1,4).Value = "PC:192.168.14.240"
$ws2.Cells.Item(1,5).Value = 1

# Row 2: 514_Teleflex -> 527_Teijin
$ws2.Cells.Item(2,1).Value = "527_Teijin"
$ws2.Cells.Item(2,2).Value = "10.101.28.176"
$ws2.Cells.Item(2,3).Value = "255.255.255.0"
$ws2.Cells.Item(2,4).Value = "XG-X2900:`t`t10.101.28.175"
$ws2.Cells.Item(2,5).Value = 1

# Row 3: 474 B_Austin -> 503_Witte
$ws2.Cells.Item(3,1).Value = "503_Witte"
$ws2.Cells.Item(3,2).Value = "192.168.0.240"
$ws2.Cells.Item(3,3).Value = "255.255.255.0"
$ws2.Cells.Item(3,4).Value = "PC:`t10.96.205.175`nNAS:`t10.96.205.166`nFH:`t10.96.205.154`n`t10.96.205.267`n-----------------------------------------`nuser:JHV_Vision, omron `nPass:*Jhv2708`n---------------------------------------`nFortiClient Austin: `nPass:`n1Pm#J@PFIkzM&Q@i `nUVt1@Ex2p78kxp30atD7we@!qGK"
$ws2.Cells.Item(3,5).Value = 1

# Row 4 (527_Teijin) got folded into row 2 above - drop the now-duplicate row
$ws2.Range("A4:E4").EntireRow.Delete()

# ----------------------------------------------------------------------
# 3) Settings
# ----------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Settings")
$ws4.Cells.Item(1,2).Value = 0
$ws4.Cells.Item(3,2).Value = 0
